$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-03 Wednesday" "2024-01-04 Thursday"

Replace-Text "72÷2=36, 0" "95÷2=47, 1"
Replace-Text "23÷3=7, 2" "10÷2=5, 0"
Replace-Text "18÷7=2, 4" "61÷8=7, 5"
Replace-Text "35÷4=8, 3" "89÷3=29, 2"
Replace-Text "81÷2=40, 1" "92÷7=13, 1"
Replace-Text "75÷8=9, 3" "26÷5=5, 1"
Replace-Text "19÷9=2, 1" "39÷2=19, 1"
Replace-Text "10÷6=1, 4" "68÷7=9, 5"
Replace-Text "65÷9=7, 2" "64÷4=16, 0"
Replace-Text "90÷9=10, 0" "98÷3=32, 2"
Replace-Text "77÷8=9, 5" "61÷3=20, 1"
Replace-Text "11÷3=3, 2" "66÷9=7, 3"
Replace-Text "71÷3=23, 2" "39÷6=6, 3"
Replace-Text "60÷4=15, 0" "42÷4=10, 2"
Replace-Text "50÷6=8, 2" "99÷5=19, 4"
Replace-Text "61÷2=30, 1" "94÷3=31, 1"
Replace-Text "33÷2=16, 1" "59÷2=29, 1"
Replace-Text "26÷9=2, 8" "81÷5=16, 1"
Replace-Text "45÷6=7, 3" "27÷8=3, 3"
Replace-Text "79÷7=11, 2" "58÷5=11, 3"
Replace-Text "79÷8=9, 7" "14÷9=1, 5"
Replace-Text "39÷4=9, 3" "61÷7=8, 5"
Replace-Text "36÷2=18, 0" "68÷7=9, 5"
Replace-Text "14÷3=4, 2" "78÷9=8, 6"
Replace-Text "61÷4=15, 1" "80÷2=40, 0"

Write-Output "Done"
